$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 2788505.39
$ws.Range("C9").Value = 434231.99
$ws.Range("D9").Value = 3222737.38
$ws.Range("E9").Value = 13.47401102847543
$ws.Range("F9").Value = 86.52598897152458
$ws.Range("G9").Value = -58.03364334702056
$ws.Range("H9").Value = -49.6434470392154
$ws.Range("I9").Value = -50.96437599356107
$ws.Range("J9").Value = 27641
$ws.Range("K9").Value = 1177
$ws.Range("L9").Value = 28818
